$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 138.36363
$ws.Cells.Item(9, 9).Value = 144.2
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 11).Value = 144.2
$ws.Cells.Item(9, 12).Value = 80
$ws.Cells.Item(9, 13).Value = 24.80000000000001
$ws.Cells.Item(9, 14).Value = -418

$ws.Cells.Item(11, 8).Value = 101.22222
$ws.Cells.Item(11, 9).Value = 101.22222
$ws.Cells.Item(11, 11).Value = 101.22222
$ws.Cells.Item(11, 13).Value = 38.77778000000001

$ws.Cells.Item(15, 8).Value = 60.5
$ws.Cells.Item(15, 9).Value = 60.5
$ws.Cells.Item(15, 11).Value = 181.5
$ws.Cells.Item(15, 13).Value = -12.5

$ws.Cells.Item(121, 8).Value = 595.6923
$ws.Cells.Item(121, 10).Value = 595.6923
$ws.Cells.Item(121, 12).Value = 1787.0769
$ws.Cells.Item(121, 14).Value = -5281.0769

$ws.Cells.Item(123, 8).Value = 41835
$ws.Cells.Item(123, 10).Value = 41835
$ws.Cells.Item(123, 12).Value = 41835
$ws.Cells.Item(123, 14).Value = -51635

$ws.Cells.Item(129, 8).Value = 1250.9753
$ws.Cells.Item(129, 9).Value = 425
$ws.Cells.Item(129, 10).Value = 1305.3158
$ws.Cells.Item(129, 11).Value = 1275
$ws.Cells.Item(129, 12).Value = 3915.9474
$ws.Cells.Item(129, 13).Value = 3725
$ws.Cells.Item(129, 14).Value = -13915.9474

$ws.Cells.Item(138, 8).Value = 2652.953
$ws.Cells.Item(138, 9).Value = 1247.3871
$ws.Cells.Item(138, 10).Value = 3459.8518
$ws.Cells.Item(138, 11).Value = 3742.1613
$ws.Cells.Item(138, 12).Value = 10379.5554
$ws.Cells.Item(138, 13).Value = 1397.8387
$ws.Cells.Item(138, 14).Value = -20659.5554

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3507.353
$ws.Cells.Item(32, 9).Value = 3469.9355
$ws.Cells.Item(32, 11).Value = 3469.9355
$ws.Cells.Item(32, 13).Value = -3182.9355

$ws.Cells.Item(61, 8).Value = 1182.963
$ws.Cells.Item(61, 9).Value = 1185
$ws.Cells.Item(61, 10).Value = 1166.6666
$ws.Cells.Item(61, 11).Value = 1185
$ws.Cells.Item(61, 12).Value = 1166.6666
$ws.Cells.Item(61, 13).Value = -973
$ws.Cells.Item(61, 14).Value = -1590.6666

$ws.Cells.Item(136, 8).Value = 1182.963
$ws.Cells.Item(136, 9).Value = 1185
$ws.Cells.Item(136, 10).Value = 1166.6666
$ws.Cells.Item(136, 11).Value = 3555
$ws.Cells.Item(136, 12).Value = 3499.9998
$ws.Cells.Item(136, 13).Value = -1005
$ws.Cells.Item(136, 14).Value = -8599.9998

$ws.Cells.Item(137, 8).Value = 41544
$ws.Cells.Item(137, 10).Value = 41544
$ws.Cells.Item(137, 12).Value = 41544
$ws.Cells.Item(137, 14).Value = -51744

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2587.2144
$ws.Cells.Item(99, 9).Value = 1293.4445
$ws.Cells.Item(99, 11).Value = 1293.4445
$ws.Cells.Item(99, 13).Value = 204.5554999999999

$ws.Cells.Item(107, 8).Value = 1267.6923
$ws.Cells.Item(107, 9).Value = 1201.25
$ws.Cells.Item(107, 10).Value = 1489.1666
$ws.Cells.Item(107, 11).Value = 1201.25
$ws.Cells.Item(107, 12).Value = 1489.1666
$ws.Cells.Item(107, 13).Value = 718.75
$ws.Cells.Item(107, 14).Value = -5329.1666

$ws.Cells.Item(134, 8).Value = 3230.8728
$ws.Cells.Item(134, 9).Value = 1079.5483
$ws.Cells.Item(134, 10).Value = 6009.6665
$ws.Cells.Item(134, 11).Value = 3238.6449
$ws.Cells.Item(134, 12).Value = 18028.9995
$ws.Cells.Item(134, 13).Value = -703.6448999999998
$ws.Cells.Item(134, 14).Value = -23098.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 206489.95
$ws.Cells.Item(31, 9).Value = 410295.72
$ws.Cells.Item(31, 10).Value = 2684.182
$ws.Cells.Item(31, 11).Value = 410295.72
$ws.Cells.Item(31, 12).Value = 2684.182
$ws.Cells.Item(31, 13).Value = -410000.72
$ws.Cells.Item(31, 14).Value = -3274.182

$ws.Cells.Item(34, 8).Value = 206489.95
$ws.Cells.Item(34, 9).Value = 410295.72
$ws.Cells.Item(34, 10).Value = 2684.182
$ws.Cells.Item(34, 11).Value = 410295.72
$ws.Cells.Item(34, 12).Value = 2684.182
$ws.Cells.Item(34, 13).Value = -410093.72
$ws.Cells.Item(34, 14).Value = -3088.182

$ws.Cells.Item(58, 8).Value = 2652.2856
$ws.Cells.Item(58, 9).Value = 1476.0714
$ws.Cells.Item(58, 10).Value = 7357.143
$ws.Cells.Item(58, 11).Value = 1476.0714
$ws.Cells.Item(58, 12).Value = 7357.143
$ws.Cells.Item(58, 13).Value = -1273.0714
$ws.Cells.Item(58, 14).Value = -7763.143

$ws.Cells.Item(132, 8).Value = 3951.9312
$ws.Cells.Item(132, 9).Value = 3185.45
$ws.Cells.Item(132, 10).Value = 5655.222
$ws.Cells.Item(132, 11).Value = 9556.349999999999
$ws.Cells.Item(132, 12).Value = 16965.666
$ws.Cells.Item(132, 13).Value = -7026.349999999999
$ws.Cells.Item(132, 14).Value = -22025.666

$ws.Cells.Item(136, 8).Value = 2652.2856
$ws.Cells.Item(136, 9).Value = 1476.0714
$ws.Cells.Item(136, 10).Value = 7357.143
$ws.Cells.Item(136, 11).Value = 4428.2142
$ws.Cells.Item(136, 12).Value = 22071.429
$ws.Cells.Item(136, 13).Value = -1878.2142
$ws.Cells.Item(136, 14).Value = -27171.429

$ws.Cells.Item(139, 8).Value = 39000
$ws.Cells.Item(139, 10).Value = 39000
$ws.Cells.Item(139, 12).Value = 39000
$ws.Cells.Item(139, 14).Value = -49280

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 4400
$ws.Cells.Item(19, 10).Value = 4400
$ws.Cells.Item(19, 12).Value = 13200
$ws.Cells.Item(19, 14).Value = -13548

$ws.Cells.Item(68, 8).Value = 1327.5278
$ws.Cells.Item(68, 9).Value = 900.4
$ws.Cells.Item(68, 10).Value = 1632.619
$ws.Cells.Item(68, 11).Value = 2701.2
$ws.Cells.Item(68, 12).Value = 4897.857
$ws.Cells.Item(68, 13).Value = -1890.2
$ws.Cells.Item(68, 14).Value = -6519.857

$ws.Cells.Item(71, 8).Value = 1327.5278
$ws.Cells.Item(71, 9).Value = 900.4
$ws.Cells.Item(71, 10).Value = 1632.619
$ws.Cells.Item(71, 11).Value = 8103.599999999999
$ws.Cells.Item(71, 12).Value = 14693.571
$ws.Cells.Item(71, 13).Value = -4047.599999999999
$ws.Cells.Item(71, 14).Value = -22805.571

$ws.Cells.Item(86, 8).Value = 1280.1333
$ws.Cells.Item(86, 9).Value = 828.5714
$ws.Cells.Item(86, 10).Value = 1675.25
$ws.Cells.Item(86, 11).Value = 2485.7142
$ws.Cells.Item(86, 12).Value = 5025.75
$ws.Cells.Item(86, 13).Value = -1299.7142
$ws.Cells.Item(86, 14).Value = -7397.75

$ws.Cells.Item(89, 8).Value = 1280.1333
$ws.Cells.Item(89, 9).Value = 828.5714
$ws.Cells.Item(89, 10).Value = 1675.25
$ws.Cells.Item(89, 11).Value = 7457.1426
$ws.Cells.Item(89, 12).Value = 15077.25
$ws.Cells.Item(89, 13).Value = -1529.1426
$ws.Cells.Item(89, 14).Value = -26933.25

$ws.Cells.Item(101, 8).Value = 5000
$ws.Cells.Item(101, 10).Value = 5000
$ws.Cells.Item(101, 12).Value = 15000
$ws.Cells.Item(101, 14).Value = -19868

$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 13).ClearContents()

$ws.Cells.Item(119, 8).Value = 1864.5
$ws.Cells.Item(119, 9).Value = 1864.5
$ws.Cells.Item(119, 11).Value = 5593.5
$ws.Cells.Item(119, 13).Value = -755.5

$ws.Cells.Item(120, 8).Value = 6515
$ws.Cells.Item(120, 9).Value = 3030
$ws.Cells.Item(120, 10).Value = 10000
$ws.Cells.Item(120, 11).Value = 9090
$ws.Cells.Item(120, 12).Value = 30000
$ws.Cells.Item(120, 13).Value = -4252
$ws.Cells.Item(120, 14).Value = -39676

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6554.718
$ws.Cells.Item(70, 9).Value = 5875.2
$ws.Cells.Item(70, 10).Value = 8819.777
$ws.Cells.Item(70, 11).Value = 5875.2
$ws.Cells.Item(70, 12).Value = 8819.777
$ws.Cells.Item(70, 13).Value = -5605.2
$ws.Cells.Item(70, 14).Value = -9359.777

$ws.Cells.Item(73, 8).Value = 6554.718
$ws.Cells.Item(73, 9).Value = 5875.2
$ws.Cells.Item(73, 10).Value = 8819.777
$ws.Cells.Item(73, 11).Value = 5875.2
$ws.Cells.Item(73, 12).Value = 8819.777
$ws.Cells.Item(73, 13).Value = -4939.2
$ws.Cells.Item(73, 14).Value = -10691.777

$ws.Cells.Item(126, 8).Value = 3401.71
$ws.Cells.Item(126, 9).Value = 2766.2112
$ws.Cells.Item(126, 10).Value = 4957.5864
$ws.Cells.Item(126, 11).Value = 8298.633600000001
$ws.Cells.Item(126, 12).Value = 14872.7592
$ws.Cells.Item(126, 13).Value = -5828.633600000001
$ws.Cells.Item(126, 14).Value = -19812.7592

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2943.0667
$ws.Cells.Item(7, 9).Value = 1328.3529
$ws.Cells.Item(7, 11).Value = 1328.3529
$ws.Cells.Item(7, 13).Value = -1216.3529

$ws.Cells.Item(40, 8).Value = 4995.8096
$ws.Cells.Item(40, 9).Value = 5193.231
$ws.Cells.Item(40, 10).Value = 4675
$ws.Cells.Item(40, 11).Value = 5193.231
$ws.Cells.Item(40, 12).Value = 4675
$ws.Cells.Item(40, 13).Value = -5057.231
$ws.Cells.Item(40, 14).Value = -4947

$ws.Cells.Item(110, 8).Value = 28800
$ws.Cells.Item(110, 10).Value = 28800
$ws.Cells.Item(110, 12).Value = 28800
$ws.Cells.Item(110, 14).Value = -36980

$ws.Cells.Item(126, 8).Value = 2943.0667
$ws.Cells.Item(126, 9).Value = 1328.3529
$ws.Cells.Item(126, 11).Value = 3985.0587
$ws.Cells.Item(126, 13).Value = -1515.0587

$ws.Cells.Item(136, 8).Value = 2461.3936
$ws.Cells.Item(136, 9).Value = 1212.3715
$ws.Cells.Item(136, 10).Value = 4142.769
$ws.Cells.Item(136, 11).Value = 3637.1145
$ws.Cells.Item(136, 12).Value = 12428.307
$ws.Cells.Item(136, 13).Value = -1087.1145
$ws.Cells.Item(136, 14).Value = -17528.307
